$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new customer record ("Joselu") was added as row 3 of the export,
# pushing every existing record (old rows 3-27) down by one (new rows 4-28).
$ws.Rows.Item(3).Insert()

# Column A (DNI_Cliente) holds digit-only IDs elsewhere in the sheet as
# plain text (e.g. "0344354326"), so force text formatting first or the
# leading zero would be lost to numeric auto-conversion.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "021247781"
$ws.Range("B3").Value = "Joselu"
$ws.Range("C3").Value = "OEOEOE"
$ws.Range("D3").Value = 91919191
$ws.Range("E3").Value = "C/ 123"
# F3 (Password) is left blank, matching the source record's empty password field.
$ws.Range("G3").Value = "aaa@eeee.net"
